$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brown, et al.")
$ws.Activate()

# B5 and B6 get the new "Sigma Aldrich" shared string
$ws.Range("B5").Value = "Sigma Aldrich"
$ws.Range("B6").Value = "Sigma Aldrich"
